$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Bump the cached "datetimeFigureOut" date field from 2/16/18 to
#    2/27/18 everywhere it is cached: the slide master's Date
#    Placeholder and every slide layout's Date Placeholder.
# ---------------------------------------------------------------------
$oldDate = "2/16/18"
$newDate = "2/27/18"

function Update-DatePlaceholder($shape) {
    if ($shape.Name -like "Date Placeholder*" -and $shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        $tr.Text = $newDate
    }
}

# Slide master.
for ($i = 1; $i -le $p.SlideMaster.Shapes.Count; $i++) {
    Update-DatePlaceholder $p.SlideMaster.Shapes.Item($i)
}

# Every slide layout hanging off the master.
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $cl = $p.SlideMaster.CustomLayouts.Item($li)
    for ($si = 1; $si -le $cl.Shapes.Count; $si++) {
        Update-DatePlaceholder $cl.Shapes.Item($si)
    }
}

# ---------------------------------------------------------------------
# 2. Slide 6 ("ELB Hands-on"): swap the webserver-failover AMI id and
#    collapse the sentence back into a single run.
# ---------------------------------------------------------------------
$oldAmi = "Launch two instances running webservers (using webserver-failover AMI ami-3ea13f29)"
$newAmi = "Launch two instances running webservers (using webserver-failover AMI ami-598b6124)"

$slide = $p.Slides.Item(6)
$content = $slide.Shapes.Item(2)
$tr = $content.TextFrame.TextRange
$fullText = $tr.Text
$start = $fullText.IndexOf($oldAmi) + 1
if ($start -gt 0) {
    $sub = $tr.Characters($start, $oldAmi.Length)
    $sub.Text = $newAmi
}
